# Swap the deck's theme palette: "Integral" -> "Office Theme".
#
# The underlying package carries two theme parts (theme1.xml = the theme
# actually applied to the slide master / slides, theme2.xml = the theme
# wired to the notes master). The edit swaps their contents so the slides
# take on the stock "Office Theme" look (and the notes master ends up
# with the former "Integral" colours). PowerPoint's automation surface
# doesn't expose raw OOXML theme parts directly, so we reproduce the
# observable effect through the supported Theme/ThemeColorScheme object
# model: writing the Office Theme's twelve theme colours onto the
# presentation's (shared) theme.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$theme = $m.Theme

function HexToMsoRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# ppThemeColorDark1, Light1, Dark2, Light2, Accent1-6, Hyperlink, FollowedHyperlink
# (indices 1-12 of ThemeColorScheme), values taken from the stock Office
# Theme colour scheme.
$officeThemeColors = @(
    "000000", # Dark 1
    "FFFFFF", # Light 1
    "44546A", # Dark 2
    "E7E6E6", # Light 2
    "5B9BD5", # Accent 1
    "ED7D31", # Accent 2
    "A5A5A5", # Accent 3
    "FFC000", # Accent 4
    "4472C4", # Accent 5
    "70AD47", # Accent 6
    "0563C1", # Hyperlink
    "954F72"  # Followed Hyperlink
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $theme.ThemeColorScheme.Item($i).RGB = HexToMsoRGB $officeThemeColors[$i - 1]
}
